# Add two new header columns (AssetName / ProjectName) and the server-side
# rows that came with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: C1 / D1, matching the bold/border/centered style of A1:B1 ---
$ws.Range("C1").Value = "AssetName"
$ws.Range("D1").Value = "ProjectName"

$ws.Range("A1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# --- New data rows 4-6 ---
$ws.Range("C4").Value = "Demo Project"
$ws.Range("D4").Value = "Demo Project"

$ws.Range("C5").Value = "Demo Project"
$ws.Range("D5").Value = "Demo Project"

$ws.Range("C6").Value = "ArcelorMittal Poland SA"
$ws.Range("D6").Value = "Arcelo"

# --- Cells that exist (so the sheet dimension grows) but hold empty text ---
# A leading apostrophe forces Excel to create a real (empty) text cell
# instead of silently clearing it; re-pasting the formats of an existing
# plain (unstyled) cell afterwards strips the quote-prefix formatting it
# leaves behind, restoring the default style.
$emptyCells = "C2", "D2", "C3", "D3", "A4", "B4", "A5", "B5", "A6", "B6"
foreach ($addr in $emptyCells) {
    $ws.Range($addr).Value = "'"
}

$ws.Range("A2").Copy()
foreach ($addr in $emptyCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
